# --- Footer / update note block rework ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New style for B22: font matching B3 (Montserrat / fontId 1) with
# horizontal=left, vertical=top alignment. Copy formats from B3 (which
# already uses that font) and then force horizontal alignment to left.
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("B22").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B22").VerticalAlignment = -4160    # xlTop

$ws.Range("B22").Value = "Actualización: mayo 2024."

# Remove the old "Ultima actualización" value that lived in D22 entirely
# (cell + formatting disappear).
$ws.Range("D22").Clear()

# Add the source/footnote text that used to live in B22 into the new B23
# cell (plain/default formatting, no explicit style).
$ws.Range("B23").Value = "Fuente: ARTF. Agencia Reguladora del Transporte Ferroviario."

# D23 keeps its existing style (s="4") but no longer carries the
# "Dirección General de Planeación" text.
$ws.Range("D23").ClearContents()

# --- Text corrections (accent fixes) ---
$ws.Range("C13").Value = "Vía Corta Oaxaca c"
$ws.Range("C19").Value = "Líneas remanentes"
